$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.146.00'
$ws.Range('D2').ClearFormats()

$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.47%  '
$ws.Range('E2').ClearFormats()

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.838.17'
$ws.Range('D3').ClearFormats()

$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.44%  '
$ws.Range('E3').ClearFormats()

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.46%  '
$ws.Range('E4').ClearFormats()

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.13'
$ws.Range('D5').ClearFormats()

$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.10%  '
$ws.Range('E5').ClearFormats()

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6182'
$ws.Range('D6').ClearFormats()

$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -2.18%  '
$ws.Range('E6').ClearFormats()

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.52%  '
$ws.Range('E7').ClearFormats()

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07486'
$ws.Range('D8').ClearFormats()

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.46%  '
$ws.Range('E8').ClearFormats()

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2927'
$ws.Range('D9').ClearFormats()

$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.39%  '
$ws.Range('E9').ClearFormats()

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.20'
$ws.Range('D10').ClearFormats()

$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +1.40%  '
$ws.Range('E10').ClearFormats()

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07709'
$ws.Range('D11').ClearFormats()

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.25%  '
$ws.Range('E11').ClearFormats()

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.841.70'
$ws.Range('D12').ClearFormats()

$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.33%  '
$ws.Range('E12').ClearFormats()

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.006'
$ws.Range('D13').ClearFormats()

$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.23%  '
$ws.Range('E13').ClearFormats()

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6727'
$ws.Range('D14').ClearFormats()

$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.27%  '
$ws.Range('E14').ClearFormats()

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '82.80'
$ws.Range('D15').ClearFormats()

$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.37%  '
$ws.Range('E15').ClearFormats()

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000009299'
$ws.Range('D16').ClearFormats()

$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -4.11%  '
$ws.Range('E16').ClearFormats()

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.948'
$ws.Range('D17').ClearFormats()

$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.16%  '
$ws.Range('E17').ClearFormats()

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '29.151.23'
$ws.Range('D18').ClearFormats()

$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.40%  '
$ws.Range('E18').ClearFormats()

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.108.89'
$ws.Range('D19').ClearFormats()

$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.04%  '
$ws.Range('E19').ClearFormats()

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '231.73'
$ws.Range('D20').ClearFormats()

$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +2.23%  '
$ws.Range('E20').ClearFormats()

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.66'
$ws.Range('D21').ClearFormats()

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.92%  '
$ws.Range('E21').ClearFormats()

$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.65%  '
$ws.Range('E22').ClearFormats()

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.178'
$ws.Range('D23').ClearFormats()

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('E23').ClearFormats()

$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.50%  '
$ws.Range('E24').ClearFormats()

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '160.32'
$ws.Range('D25').ClearFormats()

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.35%  '
$ws.Range('E25').ClearFormats()

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.534'
$ws.Range('D26').ClearFormats()

$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('E26').ClearFormats()

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1389'
$ws.Range('D27').ClearFormats()

$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.29%  '
$ws.Range('E27').ClearFormats()

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.86'
$ws.Range('D28').ClearFormats()

$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.26%  '
$ws.Range('E28').ClearFormats()

$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.58%  '
$ws.Range('E29').ClearFormats()

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.171'
$ws.Range('D30').ClearFormats()

$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.22%  '
$ws.Range('E30').ClearFormats()

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.132'
$ws.Range('D31').ClearFormats()

$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +1.43%  '
$ws.Range('E31').ClearFormats()

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.05545'
$ws.Range('D32').ClearFormats()

$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +3.34%  '
$ws.Range('E32').ClearFormats()

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.203'
$ws.Range('D33').ClearFormats()

$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.22%  '
$ws.Range('E33').ClearFormats()

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7475'
$ws.Range('D34').ClearFormats()

$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.55%  '
$ws.Range('E34').ClearFormats()

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.845'
$ws.Range('D35').ClearFormats()

$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.75%  '
$ws.Range('E35').ClearFormats()

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.142'
$ws.Range('D36').ClearFormats()

$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.25%  '
$ws.Range('E36').ClearFormats()

$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.51%  '
$ws.Range('E37').ClearFormats()

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.774'
$ws.Range('D38').ClearFormats()

$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.76%  '
$ws.Range('E38').ClearFormats()

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.225.00'
$ws.Range('D39').ClearFormats()

$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.65%  '
$ws.Range('E39').ClearFormats()

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01782'
$ws.Range('D40').ClearFormats()

$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.18%  '
$ws.Range('E40').ClearFormats()

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.489'
$ws.Range('D41').ClearFormats()

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -1.62%  '
$ws.Range('E41').ClearFormats()

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8969'
$ws.Range('D42').ClearFormats()

$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.64%  '
$ws.Range('E42').ClearFormats()

$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.51%  '
$ws.Range('E43').ClearFormats()

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.017.86'
$ws.Range('D44').ClearFormats()

$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +1.60%  '
$ws.Range('E44').ClearFormats()

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '101.97'
$ws.Range('D45').ClearFormats()

$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.45%  '
$ws.Range('E45').ClearFormats()

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '65.89'
$ws.Range('D46').ClearFormats()

$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.57%  '
$ws.Range('E46').ClearFormats()

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.00000000122'
$ws.Range('D47').ClearFormats()

$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.46%  '
$ws.Range('E47').ClearFormats()

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5108'
$ws.Range('D48').ClearFormats()

$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.11%  '
$ws.Range('E48').ClearFormats()

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4078'
$ws.Range('D49').ClearFormats()

$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.15%  '
$ws.Range('E49').ClearFormats()

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.138'
$ws.Range('D50').ClearFormats()

$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.62%  '
$ws.Range('E50').ClearFormats()

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05842'
$ws.Range('D51').ClearFormats()

$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.33%  '
$ws.Range('E51').ClearFormats()
